# Update column G ("K") values per the regenerated save_data
# (commit: "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 1
    3  = 0
    4  = 2
    5  = 2
    6  = 1
    9  = 0
    10 = 1
    11 = 1
    12 = 2
    13 = 1
    14 = 1
    15 = 0
    16 = 1
    17 = 0
    19 = 0
    20 = 1
    21 = 0
    22 = 3
    23 = 1
    24 = 1
    26 = 1
    27 = 0
    28 = 0
    29 = 1
    30 = 0
    31 = 0
    32 = 1
    33 = 2
    35 = 1
    36 = 0
    37 = 1
    38 = 0
    39 = 1
    40 = 0
    41 = 2
    42 = 0
    43 = 0
    44 = 1
    45 = 1
    46 = 1
    47 = 1
    48 = 3
    49 = 1
    50 = 0
    51 = 2
    52 = 0
    53 = 2
    54 = 0
    55 = 0
    57 = 1
    58 = 2
    59 = 2
    60 = 0
    61 = 1
    62 = 3
    63 = 1
    64 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newValues[$row]
}
